# Commit: "bezig met bekijk pagina fotos"
#
# - Add a new sheet "week 10" (cloned from the "week 7" logboek template),
#   positioned right before "Totaal".
# - Fill in week 10's first log entry (date, end time, activity text) and
#   leave the remaining rows blank, matching a fresh week's logboek.
# - Make "week 10" the active/selected tab (the previous active tab was
#   "week 7").

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("week 7")
$totaal = $wb.Worksheets.Item("Totaal")

# Clone the "week 7" sheet and drop the copy right before "Totaal".
$template.Copy($totaal)

# Excel names the clone "week 7 (2)"; rename it to the new week number.
$newSheet = $wb.Worksheets.Item("week 7 (2)")
$newSheet.Name = "week 10"

# First entry of the new week: Tuesday 2014-04-03, 08:45 - 09:05.
$newSheet.Range("B7").Value = "4/3/2014"
$newSheet.Range("D7").Value = 0.37847222222222227
$newSheet.Range("F7").Value = "bezig met fout eruit te halen"

# The remaining copied rows (2nd/3rd/4th activity of the old week) are not
# filled in yet for the new week, so clear them back out.
$newSheet.Range("C8").ClearContents()
$newSheet.Range("D8").ClearContents()
$newSheet.Range("F8").ClearContents()

$newSheet.Range("C9").ClearContents()
$newSheet.Range("D9").ClearContents()
$newSheet.Range("F9").ClearContents()

$newSheet.Range("A10").ClearContents()
$newSheet.Range("B10").ClearContents()
$newSheet.Range("C10").ClearContents()
$newSheet.Range("D10").ClearContents()
$newSheet.Range("F10").ClearContents()

# Make the new week the active sheet/selection.
$newSheet.Activate()
$newSheet.Range("F7").Select()
